$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 7859.25
$ws.Range("J69").Value = 8624.857
$ws.Range("L69").Value = 25874.571
$ws.Range("N69").Value = -27622.571

$ws.Range("H72").Value = 7859.25
$ws.Range("J72").Value = 8624.857
$ws.Range("L72").Value = 77623.713
$ws.Range("N72").Value = -86359.713

$ws.Range("H86").Value = 8705765
$ws.Range("I86").Value = 6443
$ws.Range("K86").Value = 6443
$ws.Range("M86").Value = -5320

$ws.Range("H89").Value = 8705765
$ws.Range("I89").Value = 6443
$ws.Range("K89").Value = 32215
$ws.Range("M89").Value = -26599

$ws.Range("H100").Value = 5173.294
$ws.Range("I100").Value = 5245.5
$ws.Range("K100").Value = 5245.5
$ws.Range("M100").Value = -4704.5

$ws.Range("H116").Value = 64939630
$ws.Range("J116").Value = 66689536
$ws.Range("L116").Value = 66689536
$ws.Range("N116").Value = -66696420

$ws.Range("H137").Value = 2702
$ws.Range("I137").Value = 1758.6666
$ws.Range("J137").Value = 3137.3845
$ws.Range("K137").Value = 5275.9998
$ws.Range("L137").Value = 9412.1535
$ws.Range("M137").Value = -2725.9998
$ws.Range("N137").Value = -14512.1535

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2219.8
$ws.Range("I45").Value = 2166.6667
$ws.Range("J45").Value = 2299.5
$ws.Range("K45").Value = 2166.6667
$ws.Range("L45").Value = 2299.5
$ws.Range("M45").Value = -1789.6667
$ws.Range("N45").Value = -3053.5

$ws.Range("H132").Value = 30305492
$ws.Range("I132").Value = 38463850
$ws.Range("K132").Value = 115391550
$ws.Range("M132").Value = -115389020

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1012.8571
$ws.Range("I99").Value = 975.38464
$ws.Range("K99").Value = 975.38464
$ws.Range("M99").Value = 522.61536

$ws.Range("H134").Value = 2998.5
$ws.Range("I134").Value = 1568.8572
$ws.Range("J134").Value = 5000
$ws.Range("K134").Value = 4706.571599999999
$ws.Range("L134").Value = 15000
$ws.Range("M134").Value = -2171.571599999999
$ws.Range("N134").Value = -20070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1670.4706
$ws.Range("I122").Value = 1681.125
$ws.Range("J122").Value = 1500
$ws.Range("K122").Value = 5043.375
$ws.Range("L122").Value = 4500
$ws.Range("M122").Value = -2593.375
$ws.Range("N122").Value = -9400

$ws.Range("H134").Value = 1917.7576
$ws.Range("I134").Value = 1235.25
$ws.Range("J134").Value = 3737.7778
$ws.Range("K134").Value = 3705.75
$ws.Range("L134").Value = 11213.3334
$ws.Range("M134").Value = -1170.75
$ws.Range("N134").Value = -16283.3334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 57.333332
$ws.Range("J33").Value = 80
$ws.Range("L33").Value = 480
$ws.Range("N33").Value = -1046

$ws.Range("H34").Value = 2209.9375
$ws.Range("I34").Value = 679
$ws.Range("J34").Value = 2905.818
$ws.Range("K34").Value = 2037
$ws.Range("L34").Value = 8717.454000000002
$ws.Range("M34").Value = -1953
$ws.Range("N34").Value = -8885.454000000002

$ws.Range("H39").Value = 2938.8572
$ws.Range("J39").Value = 3906.4
$ws.Range("L39").Value = 11719.2
$ws.Range("N39").Value = -12307.2

$ws.Range("H96").Value = 10028
$ws.Range("J96").Value = 10028
$ws.Range("L96").Value = 30084
$ws.Range("N96").Value = -34202

$ws.Range("H99").Value = 650
$ws.Range("I99").Value = 650
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 1950
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = $null
$ws.Range("N99").Value = 296

$ws.Range("H103").Value = 870
$ws.Range("I103").Value = 202.5
$ws.Range("J103").Value = 1537.5
$ws.Range("K103").Value = 607.5
$ws.Range("L103").Value = 4612.5
$ws.Range("M103").Value = 271.5
$ws.Range("N103").Value = -6370.5

$ws.Range("H104").Value = 2440.3333
$ws.Range("I104").Value = 1161.5
$ws.Range("K104").Value = 3484.5
$ws.Range("M104").Value = -863.5

$ws.Range("H108").Value = 2467.5715
$ws.Range("I108").Value = 2212.1667
$ws.Range("J108").Value = 4000
$ws.Range("K108").Value = 6636.500100000001
$ws.Range("L108").Value = 12000
$ws.Range("M108").Value = -3756.500100000001
$ws.Range("N108").Value = -17760

$ws.Range("H111").Value = 916
$ws.Range("I111").Value = 916
$ws.Range("K111").Value = 2748
$ws.Range("M111").Value = 319

$ws.Range("H112").Value = 2541.8572
$ws.Range("I112").Value = 2832.1667
$ws.Range("K112").Value = 8496.500100000001
$ws.Range("M112").Value = -7388.500100000001

$ws.Range("H115").Value = 6718.85
$ws.Range("J115").Value = 6718.85
$ws.Range("L115").Value = 20156.55
$ws.Range("N115").Value = -22506.55

$ws.Range("H116").Value = 2579.8
$ws.Range("I116").Value = 299.66666
$ws.Range("K116").Value = 898.9999799999999
$ws.Range("M116").Value = 2543.00002

$ws.Range("H118").Value = 2069.4
$ws.Range("I118").Value = 449.33334
$ws.Range("J118").Value = 4499.5
$ws.Range("K118").Value = 1348.00002
$ws.Range("L118").Value = 13498.5
$ws.Range("M118").Value = -105.0000199999999
$ws.Range("N118").Value = -15984.5

$ws.Range("H120").Value = 31766.25
$ws.Range("I120").Value = 23499.5
$ws.Range("K120").Value = 70498.5
$ws.Range("M120").Value = -65660.5

$ws.Range("H121").Value = 84451.31
$ws.Range("I121").Value = 11746.111
$ws.Range("K121").Value = 35238.333
$ws.Range("M121").Value = -33928.333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1128.1111
$ws.Range("J82").Value = 984
$ws.Range("L82").Value = 984
$ws.Range("N82").Value = -1706

$ws.Range("H85").Value = 1128.1111
$ws.Range("J85").Value = 984
$ws.Range("L85").Value = 984
$ws.Range("N85").Value = -3480

$ws.Range("H132").Value = 5580.515
$ws.Range("I132").Value = 3488.2222
$ws.Range("J132").Value = 8091.2666
$ws.Range("K132").Value = 10464.6666
$ws.Range("L132").Value = 24273.7998
$ws.Range("M132").Value = -7934.6666
$ws.Range("N132").Value = -29333.7998

$ws.Range("H134").Value = 47653
$ws.Range("J134").Value = 47653
$ws.Range("L134").Value = 47653
$ws.Range("N134").Value = -57793

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 12999.833
$ws.Range("I15").Value = 9500
$ws.Range("J15").Value = 13699.8
$ws.Range("K15").Value = 9500
$ws.Range("L15").Value = 13699.8
$ws.Range("M15").Value = -9212
$ws.Range("N15").Value = -14275.8

$ws.Range("H125").Value = 50000
$ws.Range("J125").Value = 50000
$ws.Range("L125").Value = 50000
$ws.Range("N125").Value = -59840

$ws.Range("H132").Value = 2859.3215
$ws.Range("I132").Value = 2230.2046
$ws.Range("J132").Value = 5166.0835
$ws.Range("K132").Value = 6690.6138
$ws.Range("L132").Value = 15498.2505
$ws.Range("M132").Value = -4160.6138
$ws.Range("N132").Value = -20558.2505
